$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$value = 88.80208333333334

# Columns B..U correspond to epoch 50..1000 in steps of 50
for ($i = 0; $i -lt 20; $i++) {
    $col = 2 + $i
    $epoch = 50 * ($i + 1)
    $ws.Cells.Item(1, $col).Value = "epoch$epoch"
    $ws.Cells.Item(2, $col).Value = $value
}
